$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells O1, P1 with same style as existing header (copy style from N1)
$ws.Cells.Item(1,14).Copy($ws.Cells.Item(1,15))
$ws.Cells.Item(1,14).Copy($ws.Cells.Item(1,16))
$ws.Cells.Item(1,15).Value = "RDU"
$ws.Cells.Item(1,16).Value = "capacité au champs"

# Fill O and P columns (constant values) for data rows 2..97
for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r,15).Value = 0.5912372199999999
    $ws.Cells.Item($r,16).Value = 1.77371166
}

# Update KC (I), ETR (J), ETR_PV (N) columns where KC (the irrigation
# coefficient) changes. ETR = ETP * KC and ETR_PV = ETP_PV * KC, so both
# are recomputed from the unchanged ETP (H) / ETP_PV (M) columns using the
# new KC value.
for ($r = 2; $r -le 33; $r++) {
    $h = $ws.Cells.Item($r,8).Value2
    $m = $ws.Cells.Item($r,13).Value2
    $ws.Cells.Item($r,9).Value = 0.5
    $ws.Cells.Item($r,10).Value = $h * 0.5
    $ws.Cells.Item($r,14).Value = $m * 0.5
}
for ($r = 34; $r -le 76; $r++) {
    $h = $ws.Cells.Item($r,8).Value2
    $m = $ws.Cells.Item($r,13).Value2
    $ws.Cells.Item($r,9).Value = 1
    $ws.Cells.Item($r,10).Value = $h * 1
    $ws.Cells.Item($r,14).Value = $m * 1
}
# rows 77-84 already have KC = 1, so nothing to update there.
for ($r = 85; $r -le 97; $r++) {
    $h = $ws.Cells.Item($r,8).Value2
    $m = $ws.Cells.Item($r,13).Value2
    $ws.Cells.Item($r,9).Value = 0.7
    $ws.Cells.Item($r,10).Value = $h * 0.7
    $ws.Cells.Item($r,14).Value = $m * 0.7
}
